# Auto-generated Excel COM-interop script to apply the Maduin_Profits update
# Updates columns H-N (currentAveragePrice.. LeveProfitHQ) for specific rows across all 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 4014.7
$ws.Range("J40").Value = 5950
$ws.Range("L40").Value = 5950
$ws.Range("N40").Value = -6300

$ws.Range("H86").Value = 7786.143
$ws.Range("I86").Value = 6800
$ws.Range("J86").Value = 8180.6
$ws.Range("K86").Value = 6800
$ws.Range("L86").Value = 8180.6
$ws.Range("M86").Value = -5677
$ws.Range("N86").Value = -10426.6

$ws.Range("H89").Value = 7786.143
$ws.Range("I89").Value = 6800
$ws.Range("J89").Value = 8180.6
$ws.Range("K89").Value = 34000
$ws.Range("L89").Value = 40903
$ws.Range("M89").Value = -28384
$ws.Range("N89").Value = -52135

$ws.Range("H137").Value = 2353.5715


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1940.6666
$ws.Range("I2").Value = 154.75
$ws.Range("J2").Value = 5512.5
$ws.Range("K2").Value = 154.75
$ws.Range("L2").Value = 5512.5
$ws.Range("M2").Value = -41.75
$ws.Range("N2").Value = -5738.5

$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490

$ws.Range("H116").Value = 1940.6666
$ws.Range("I116").Value = 154.75
$ws.Range("J116").Value = 5512.5
$ws.Range("K116").Value = 154.75
$ws.Range("L116").Value = 5512.5
$ws.Range("M116").Value = 2139.25
$ws.Range("N116").Value = -10100.5

$ws.Range("H132").Value = 4998.3335
$ws.Range("I132").Value = 5997.5
$ws.Range("K132").Value = 17992.5
$ws.Range("M132").Value = -15462.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1940.6666
$ws.Range("I3").Value = 154.75
$ws.Range("J3").Value = 5512.5
$ws.Range("K3").Value = 154.75
$ws.Range("L3").Value = 5512.5
$ws.Range("M3").Value = -40.75
$ws.Range("N3").Value = -5740.5

$ws.Range("H20").Value = 1730.5714
$ws.Range("I20").Value = 2158.8
$ws.Range("J20").Value = 660
$ws.Range("K20").Value = 2158.8
$ws.Range("L20").Value = 660
$ws.Range("M20").Value = -1911.8
$ws.Range("N20").Value = -1154


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9000
$ws.Range("I51").Value = 9000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8264

$ws.Range("H61").Value = 9000
$ws.Range("I61").Value = 9000
$ws.Range("K61").Value = 9000
$ws.Range("M61").Value = -8652

$ws.Range("I94").Value = 3064.889
$ws.Range("J94").Value = 2951.6667
$ws.Range("K94").Value = 3064.889
$ws.Range("L94").Value = 2951.6667
$ws.Range("M94").Value = -2613.889
$ws.Range("N94").Value = -3853.6667

$ws.Range("H99").Value = 6555
$ws.Range("I99").Value = 5719.6
$ws.Range("K99").Value = 5719.6
$ws.Range("M99").Value = -4221.6

$ws.Range("H126").Value = 6555
$ws.Range("I126").Value = 5719.6
$ws.Range("K126").Value = 17158.8
$ws.Range("M126").Value = -14688.8

$ws.Range("H133").Value = 70326
$ws.Range("J133").Value = 70326
$ws.Range("L133").Value = 70326
$ws.Range("N133").Value = -75386

$ws.Range("H134").Value = 2435.4707
$ws.Range("I134").Value = 2274.3333
$ws.Range("J134").Value = 2822.2
$ws.Range("K134").Value = 6822.999899999999
$ws.Range("L134").Value = 8466.599999999999
$ws.Range("M134").Value = -4287.999899999999
$ws.Range("N134").Value = -13536.6

$ws.Range("H141").Value = 312571.9
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 312571.9
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 312571.9
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -322931.9


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1346.3334
$ws.Range("I17").Value = 192.66667
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 578.00001
$ws.Range("L17").Value = 7500
$ws.Range("M17").Value = -409.00001
$ws.Range("N17").Value = -7838

$ws.Range("H22").Value = 175000
$ws.Range("I22").Value = 250000
$ws.Range("K22").Value = 750000
$ws.Range("M22").Value = -749831

$ws.Range("H27").Value = 175000
$ws.Range("I27").Value = 250000
$ws.Range("K27").Value = 750000
$ws.Range("M27").Value = -749898

$ws.Range("H113").Value = 428.66666
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H128").Value = 499991.34
$ws.Range("I128").Value = 499991.34
$ws.Range("K128").Value = 1499974.02
$ws.Range("M128").Value = -1494994.02


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 74.36364
$ws.Range("I2").Value = 13.285714
$ws.Range("K2").Value = 13.285714
$ws.Range("M2").Value = 99.714286

$ws.Range("H70").Value = 9868.799999999999
$ws.Range("I70").Value = 10836.5
$ws.Range("K70").Value = 10836.5
$ws.Range("M70").Value = -10566.5

$ws.Range("H73").Value = 9868.799999999999
$ws.Range("I73").Value = 10836.5
$ws.Range("K73").Value = 10836.5
$ws.Range("M73").Value = -9900.5

$ws.Range("H102").Value = 1065
$ws.Range("I102").Value = 1065
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1065
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 557
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 783.8
$ws.Range("I107").Value = 534.4286
$ws.Range("J107").Value = 1365.6666
$ws.Range("K107").Value = 534.4286
$ws.Range("L107").Value = 1365.6666
$ws.Range("M107").Value = 1385.5714
$ws.Range("N107").Value = -5205.6666

$ws.Range("H113").Value = 2000.6666
$ws.Range("I113").Value = 2000.6666
$ws.Range("K113").Value = 2000.6666
$ws.Range("M113").Value = 169.3334

$ws.Range("H132").Value = 3947.0667
$ws.Range("I132").Value = 3946.8462
$ws.Range("K132").Value = 11840.5386
$ws.Range("M132").Value = -9310.5386


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H46").Value = 2376.4707
$ws.Range("I46").Value = 1733.3334
$ws.Range("J46").Value = 2727.2727
$ws.Range("K46").Value = 1733.3334
$ws.Range("L46").Value = 2727.2727
$ws.Range("M46").Value = -1545.3334
$ws.Range("N46").Value = -3103.2727

$ws.Range("H122").Value = 3307.4
$ws.Range("I122").Value = 1996.3334
$ws.Range("J122").Value = 5274
$ws.Range("K122").Value = 5989.0002
$ws.Range("L122").Value = 15822
$ws.Range("M122").Value = -3539.0002
$ws.Range("N122").Value = -20722

$ws.Range("H125").Value = 88143.60000000001
$ws.Range("J125").Value = 88143.60000000001
$ws.Range("L125").Value = 88143.60000000001
$ws.Range("N125").Value = -97983.60000000001


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H64").Value = 29995
$ws.Range("I64").Value = 29990
$ws.Range("K64").Value = 29990
$ws.Range("M64").Value = -29742

$ws.Range("H67").Value = 29995
$ws.Range("I67").Value = 29990
$ws.Range("K67").Value = 29990
$ws.Range("M67").Value = -29132

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H126").Value = 1995
$ws.Range("I126").Value = 1995
$ws.Range("J126").Value = 1995
$ws.Range("K126").Value = 5985
$ws.Range("L126").Value = 5985
$ws.Range("M126").Value = -3515
$ws.Range("N126").Value = -10925

$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

